$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 13842.286
$ws.Range("I51").Value = 7785.8335
$ws.Range("J51").Value = 16264.866
$ws.Range("K51").Value = 7785.8335
$ws.Range("L51").Value = 16264.866
$ws.Range("M51").Value = -7301.8335
$ws.Range("N51").Value = -17232.866
$ws.Range("H62").Value = 402279.8
$ws.Range("I62").Value = 801760.6
$ws.Range("K62").Value = 801760.6
$ws.Range("M62").Value = -801136.6
$ws.Range("H64").Value = 6980.8335
$ws.Range("I64").Value = 5500
$ws.Range("K64").Value = 5500
$ws.Range("M64").Value = -5252
$ws.Range("H65").Value = 402279.8
$ws.Range("I65").Value = 801760.6
$ws.Range("K65").Value = 4008803
$ws.Range("M65").Value = -4005683
$ws.Range("H67").Value = 6980.8335
$ws.Range("I67").Value = 5500
$ws.Range("K67").Value = 5500
$ws.Range("M67").Value = -4642
$ws.Range("H88").Value = 4861.143
$ws.Range("J88").Value = 4861.143
$ws.Range("L88").Value = 4861.143
$ws.Range("N88").Value = -5673.143
$ws.Range("H91").Value = 4861.143
$ws.Range("J91").Value = 4861.143
$ws.Range("L91").Value = 4861.143
$ws.Range("N91").Value = -7669.143
$ws.Range("H97").Value = 6048.5
$ws.Range("I97").Value = 1098
$ws.Range("J97").Value = 10999
$ws.Range("K97").Value = 3294
$ws.Range("L97").Value = 32997
$ws.Range("M97").Value = -2798
$ws.Range("N97").Value = -33989
$ws.Range("H107").Value = 1018.3333
$ws.Range("I107").Value = 819.2143
$ws.Range("K107").Value = 819.2143
$ws.Range("M107").Value = 1100.7857
$ws.Range("H112").Value = 3076.147
$ws.Range("J112").Value = 3149.6875
$ws.Range("L112").Value = 9449.0625
$ws.Range("N112").Value = -11665.0625
$ws.Range("H116").Value = 3440.4348
$ws.Range("I116").Value = 3304.5557
$ws.Range("K116").Value = 3304.5557
$ws.Range("M116").Value = 137.4443000000001
$ws.Range("H135").Value = 10824.083
$ws.Range("I135").Value = 2543
$ws.Range("K135").Value = 22887
$ws.Range("M135").Value = -20352
$ws.Range("H138").Value = 6310.298
$ws.Range("I138").Value = 3064.9167
$ws.Range("J138").Value = 7423
$ws.Range("K138").Value = 9194.750100000001
$ws.Range("L138").Value = 22269
$ws.Range("M138").Value = -4054.750100000001
$ws.Range("N138").Value = -32549

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = ""
$ws.Range("N52").Value = 0
$ws.Range("H74").Value = 2599.3333
$ws.Range("I74").Value = 1399
$ws.Range("J74").Value = 3199.5
$ws.Range("K74").Value = 1399
$ws.Range("L74").Value = 3199.5
$ws.Range("M74").Value = -525
$ws.Range("N74").Value = -4947.5
$ws.Range("H77").Value = 2599.3333
$ws.Range("I77").Value = 1399
$ws.Range("J77").Value = 3199.5
$ws.Range("K77").Value = 6995
$ws.Range("L77").Value = 15997.5
$ws.Range("M77").Value = -2627
$ws.Range("N77").Value = -24733.5
$ws.Range("H122").Value = 1452.8422
$ws.Range("I122").Value = 1296.2
$ws.Range("K122").Value = 3888.6
$ws.Range("M122").Value = -1438.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 16561.568
$ws.Range("J16").Value = 32761.273
$ws.Range("L16").Value = 32761.273
$ws.Range("N16").Value = -33335.273
$ws.Range("H31").Value = 5786.567
$ws.Range("I31").Value = 7622.846
$ws.Range("K31").Value = 7622.846
$ws.Range("M31").Value = -7327.846
$ws.Range("H34").Value = 5786.567
$ws.Range("I34").Value = 7622.846
$ws.Range("K34").Value = 7622.846
$ws.Range("M34").Value = -7420.846
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = ""
$ws.Range("N57").Value = 0
$ws.Range("H62").Value = 7642
$ws.Range("I62").Value = 6936.8
$ws.Range("J62").Value = 8082.75
$ws.Range("K62").Value = 6936.8
$ws.Range("L62").Value = 8082.75
$ws.Range("M62").Value = -6312.8
$ws.Range("N62").Value = -9330.75
$ws.Range("H65").Value = 7642
$ws.Range("I65").Value = 6936.8
$ws.Range("J65").Value = 8082.75
$ws.Range("K65").Value = 34684
$ws.Range("L65").Value = 40413.75
$ws.Range("M65").Value = -31564
$ws.Range("N65").Value = -46653.75
$ws.Range("H113").Value = 16561.568
$ws.Range("J113").Value = 32761.273
$ws.Range("L113").Value = 32761.273
$ws.Range("N113").Value = -37101.273
$ws.Range("H132").Value = 3666.825
$ws.Range("I132").Value = 3516.1936
$ws.Range("K132").Value = 10548.5808
$ws.Range("M132").Value = -8018.5808
$ws.Range("H134").Value = 4460.75
$ws.Range("I134").Value = 3723.1
$ws.Range("J134").Value = 8149
$ws.Range("K134").Value = 11169.3
$ws.Range("L134").Value = 24447
$ws.Range("M134").Value = -8634.299999999999
$ws.Range("N134").Value = -29517

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2654.4546
$ws.Range("I68").Value = 1826
$ws.Range("J68").Value = 2965.125
$ws.Range("K68").Value = 5478
$ws.Range("L68").Value = 8895.375
$ws.Range("M68").Value = -4667
$ws.Range("N68").Value = -10517.375
$ws.Range("H71").Value = 2654.4546
$ws.Range("I71").Value = 1826
$ws.Range("J71").Value = 2965.125
$ws.Range("K71").Value = 16434
$ws.Range("L71").Value = 26686.125
$ws.Range("M71").Value = -12378
$ws.Range("N71").Value = -34798.125
$ws.Range("H113").Value = 3857.2856
$ws.Range("J113").Value = 4383.273
$ws.Range("L113").Value = 13149.819
$ws.Range("N113").Value = -17489.819
$ws.Range("H138").Value = 8398.700000000001
$ws.Range("I138").Value = 3427
$ws.Range("K138").Value = 10281
$ws.Range("M138").Value = -5141

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 278.72223
$ws.Range("I107").Value = 267.72726
$ws.Range("J107").Value = 296
$ws.Range("K107").Value = 267.72726
$ws.Range("L107").Value = 296
$ws.Range("M107").Value = 1652.27274
$ws.Range("N107").Value = -4136
$ws.Range("H113").Value = 10167.941
$ws.Range("I113").Value = 5246.577
$ws.Range("K113").Value = 5246.577
$ws.Range("M113").Value = -3076.577
$ws.Range("H126").Value = 4745.0835
$ws.Range("I126").Value = 4464.1
$ws.Range("J126").Value = 6150
$ws.Range("K126").Value = 13392.3
$ws.Range("L126").Value = 18450
$ws.Range("M126").Value = -10922.3
$ws.Range("N126").Value = -23390

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6233.1333
$ws.Range("I46").Value = 6289.222
$ws.Range("K46").Value = 6289.222
$ws.Range("M46").Value = -6101.222
$ws.Range("H61").Value = 2740.5652
$ws.Range("I61").Value = 2878.375
$ws.Range("K61").Value = 2878.375
$ws.Range("M61").Value = -2676.375
$ws.Range("H113").Value = 2740.5652
$ws.Range("I113").Value = 2878.375
$ws.Range("K113").Value = 2878.375
$ws.Range("M113").Value = -708.375
$ws.Range("H132").Value = 22319.342
$ws.Range("I132").Value = 28768.643
$ws.Range("J132").Value = 4261.3
$ws.Range("K132").Value = 86305.929
$ws.Range("L132").Value = 12783.9
$ws.Range("M132").Value = -83775.929
$ws.Range("N132").Value = -17843.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 17661
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 17661
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = ""
$ws.Range("N108").Value = 0
$ws.Range("H122").Value = 4493.8237
$ws.Range("I122").Value = 3313.9285
$ws.Range("K122").Value = 9941.7855
$ws.Range("M122").Value = -7491.7855
$ws.Range("H126").Value = 4499.8
$ws.Range("J126").Value = 4831.6665
$ws.Range("L126").Value = 14494.9995
$ws.Range("N126").Value = -19434.9995
$ws.Range("H132").Value = 2977.3333
$ws.Range("I132").Value = 2876.2
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8628.599999999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6098.599999999999
$ws.Range("N132").Value = -20060
